# Add new columns I (I0) and J (IF) to the sheet, matching the style of
# the existing header row (column H's style) and filling in the data
# for rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy the style from the existing H1 header cell so the new
# headers look consistent (bold, bordered, centered), then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns, keyed by row number.
$values = @{
    2  = @(5, 5)
    3  = @(6, 6)
    4  = @(6, 7)
    5  = @(9, 9)
    6  = @(4, 6)
    7  = @(6, 7)
    8  = @(5, 7)
    9  = @(8, 8)
    10 = @(5, 6)
    11 = @(8, 8)
    12 = @(7, 9)
    13 = @(7, 9)
    14 = @(7, 8)
    15 = @(8, 8)
    16 = @(7, 7)
    17 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
